# Update the "Estado de Cuenta" worker-period table (rows 16-29 on Hoja1).
# The data is being reorganized: instead of being grouped by worker
# (LUIS GABRIEL BROCHERO MARTINEZ for periods 1911->1905, then FREY ALEJANDRO
# PEREZ MARTINEZ for periods 1911->1905), the rows are now grouped by period
# (1905->1911 ascending), alternating between the two workers, and the
# "Valor Mora" amount of 26500 (vs. the regular 33125) now applies to the
# most recent period (1911) instead of the oldest period shown (1911 before
# the re-sort, i.e. the first row of each worker's block).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$tipoDoc = "CC"

$trabajador1Doc    = "73210713"
$trabajador1Nombre = "LUIS GABRIEL BROCHERO MARTINEZ"

$trabajador2Doc    = "15050164"
$trabajador2Nombre = "FREY ALEJANDRO PEREZ MARTINEZ"

$salarioBasico = 828116

# Ordered periods (ascending) and the corresponding "Valor Mora" value.
# Every period uses 33125 except the most recent one (1911), which uses 26500.
$periodos = @("1905", "1906", "1907", "1908", "1909", "1910", "1911")

$row = 16
foreach ($periodo in $periodos) {
    if ($periodo -eq "1911") {
        $valorMora = 26500
    } else {
        $valorMora = 33125
    }

    # Worker 1 row
    $ws.Cells.Item($row, 2).Value = $tipoDoc
    $ws.Cells.Item($row, 3).Value = $trabajador1Doc
    $ws.Cells.Item($row, 4).Value = $trabajador1Nombre
    $ws.Cells.Item($row, 5).Value = $periodo
    $ws.Cells.Item($row, 6).Value = $valorMora
    $ws.Cells.Item($row, 7).Value = $salarioBasico
    $row = $row + 1

    # Worker 2 row
    $ws.Cells.Item($row, 2).Value = $tipoDoc
    $ws.Cells.Item($row, 3).Value = $trabajador2Doc
    $ws.Cells.Item($row, 4).Value = $trabajador2Nombre
    $ws.Cells.Item($row, 5).Value = $periodo
    $ws.Cells.Item($row, 6).Value = $valorMora
    $ws.Cells.Item($row, 7).Value = $salarioBasico
    $row = $row + 1
}
